$d = $word.ActiveDocument

# Every paragraph in the body gets an explicit pageBreakBefore="0" in its pPr.
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# The heading / title / subtitle styles (the ones with keepNext/keepLines
# already set) also get an explicit pageBreakBefore="0" in their pPr.
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    $s.ParagraphFormat.PageBreakBefore = 0
}
